# Generate Report for Handoff
# Updates the localization status workbook:
#  - Marks the df8c8319 file as "Ready for handoff" (was "In Translation")
#    on the Overview, zh-cn and de-de sheets.
#  - Refreshes the "Latest Handoff Datetime" timestamps for that file.
#  - Sets Priority to "mt" for that file's zh-cn / de-de rows.
#  - Widens the affected datetime columns to fit the new content.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 3 is the df8c8319-... file ---
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-03 08:17:13"

# --- zh-cn sheet ---
# Row 2: 9760236c-... file -> status updated, handoff datetime re-stamped (same value)
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-09-03 08:16:31"

# Row 3: df8c8319-... file -> status updated, priority changed, new handoff datetime
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "mt"
$zhcn.Range("H3").Value = "2016-09-03 08:17:09"

# --- de-de sheet ---
# Row 2: 9760236c-... file -> status updated, handoff datetime unchanged
$dede.Range("C2").Value = "Ready for handoff"

# Row 3: df8c8319-... file -> status updated, priority changed, new handoff datetime
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "mt"
$dede.Range("H3").Value = "2016-09-03 08:17:13"

# --- Column width adjustments (datetime columns grew wider to fit new text) ---
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33
$zhcn.Columns.Item(3).ColumnWidth = 16.33
$dede.Columns.Item(3).ColumnWidth = 16.33
